# Weekly fruit/vegetable price update: Comercializadora del Agro de Limari - Palta (avocado)
# Inserts 5 new price records (row 486-490, date 2022-12-19 / serial 44889) and shifts
# the remaining historical records down by 5 rows (dimension grows from T568 to T573).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44889,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Primera",360,1700,1800,1750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44889,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Segunda",300,1500,1600,1550,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1550,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44889,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",400,2400,2500,2450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44889,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",320,2100,2200,2150,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44889,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Tercera",280,1700,1800,1750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44742,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",300,3700,3800,3750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44742,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,3400,3500,3450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44742,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,2900,3000,2950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44811,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","1a nueva(o)",300,1900,2000,1950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44811,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","2a nueva(o)",300,1700,1800,1750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44811,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","3a nueva (o)",200,1200,1300,1250,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1250,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44441,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Especial",200,2450,2500,2475,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2475,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44441,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Primera",240,2250,2300,2275,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2275,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44441,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Segunda",200,2000,2100,2050,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2050,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44441,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","1a nueva(o)",400,1700,1800,1750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44441,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","2a nueva(o)",360,1500,1600,1550,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1550,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44441,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","3a nueva (o)",200,1300,1400,1350,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1350,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44244,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",200,4100,4200,4150,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44244,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",200,3800,3900,3850,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3850,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44244,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",200,3300,3400,3350,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3350,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44272,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",200,5000,5100,5050,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",5050,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44272,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",200,4800,4900,4850,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4850,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44272,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,4600,4700,4650,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4650,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44581,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",300,2500,2600,2550,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2550,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44581,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,2200,2300,2250,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2250,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44581,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,1900,2000,1950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44161,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",300,3350,3400,3375,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3375,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44161,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,3150,3200,3175,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3175,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44161,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",300,2750,2800,2775,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",2775,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44685,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",300,2900,3000,2950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44685,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,2700,2800,2750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44685,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,2400,2500,2450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44762,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",240,3900,4000,3950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44762,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,3700,3800,3750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44762,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,3400,3500,3450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44762,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Negra de La Cruz","Primera",200,1400,1500,1450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44762,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Negra de La Cruz","Segunda",160,1100,1200,1150,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44294,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",200,5250,5300,5275,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",5275,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44294,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",240,5000,5100,5050,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",5050,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44294,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,4850,4900,4875,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4875,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44413,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Especial",200,2450,2500,2475,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2475,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44413,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Primera",300,2150,2200,2175,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2175,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44413,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Segunda",200,1850,1900,1875,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1875,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44630,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",400,2900,3000,2950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44630,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,2600,2700,2650,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2650,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44630,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",300,2300,2400,2350,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2350,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44595,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",340,2500,2600,2550,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2550,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44595,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",400,2200,2300,2250,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2250,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44595,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",300,1800,1900,1850,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1850,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44203,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",200,4750,4800,4775,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4775,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44203,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",200,4450,4500,4475,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4475,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44203,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",200,4150,4200,4175,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4175,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44503,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",340,2500,2600,2550,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2550,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44503,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",400,2200,2300,2250,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2250,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44503,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",400,1900,2000,1950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44503,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Tercera",360,1500,1600,1550,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1550,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44763,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",200,3900,4000,3950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44763,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",240,3700,3800,3750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44763,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",200,3400,3500,3450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44763,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Negra de La Cruz","Primera",200,1400,1500,1450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44763,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Negra de La Cruz","Segunda",160,1100,1200,1150,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44455,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Especial",360,2100,2200,2150,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44455,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Primera",400,1900,2000,1950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44455,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Edranol","Segunda",340,1600,1700,1650,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1650,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44455,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","1a nueva(o)",400,2300,2400,2350,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2350,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44455,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","2a nueva(o)",300,2100,2200,2150,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44622,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",240,2800,2900,2850,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2850,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44622,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",400,2600,2700,2650,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2650,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44622,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",300,2200,2300,2250,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2250,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44657,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",400,2600,2700,2650,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2650,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44657,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",400,2300,2400,2350,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2350,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44657,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",300,1900,2000,1950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",1950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44636,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",300,3000,3100,3050,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3050,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44636,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,2800,2900,2850,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2850,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44636,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",200,2600,2700,2650,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2650,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44741,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",300,3700,3800,3750,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3750,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44741,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,3400,3500,3450,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3450,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44741,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,2900,3000,2950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44252,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",200,4250,4300,4275,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",4275,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44252,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",240,3950,4000,3975,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3975,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44252,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",200,3650,3700,3675,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3675,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44364,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Fuerte","Especial",240,3100,3200,3150,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3150,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44364,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Fuerte","Primera",300,2900,3000,2950,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2950,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44714,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Especial",240,3200,3300,3250,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3250,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44714,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Primera",300,3000,3100,3050,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",3050,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44714,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Hass","Segunda",240,2800,2900,2850,"`$/kilo (en caja de 17 kilos)","Provincia de Limarí",2850,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44335,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Fuerte","Especial",200,2750,2800,2775,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",2775,1)
    ,(2,"Comercializadora del Agro de Limarí","Coquimbo",44335,4,"Fruta",100106,"Oleaginosos",100106002,"Palta","Fuerte","Primera",240,2550,2600,2575,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",2575,1)
)

$startRow = 486
$nRows = $rows.Count
$nCols = 20
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$endRow = $startRow + $nRows - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 20)).Value = $arr

# New rows 569-573 are brand-new cells; make sure column D keeps the date display format
# used throughout the rest of the "Fecha" column.
$ws.Range("D569:D573").NumberFormat = "YYYY-MM-DD HH:MM:SS"
